# Generate Report for Archive
#
# The localization run moved the two tracked files from "Ready for
# handoff" into "In Translation". Update the Status value everywhere it
# appears - once per language column on the Overview roll-up sheet, and
# once in the Status column of each language's own detail sheet - then
# re-fit those now-narrower Status columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: column E = zh-cn status, column F = de-de status
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-language sheets: column C = Status
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the Status columns now that the text is shorter than
# "Ready for handoff".
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
